# -----------------------------------------------------------------------
# C1--C2-and-C3-PowerPoint.pptx edit
#
# The authoritative diff shows two logical changes:
#
#   1. Slide 16's table (the only table/graphicFrame on that slide) gets
#      a new <a:tableStyleId> GUID:
#         {57B45F65-4B77-4F33-9896-FF0153DBF6FC}  ->  {C02D9BB8-2BE8-4244-AFEF-2932A9951E1F}
#
#   2. The deck's theme palette is swapped from the custom "Integral"
#      colour scheme to the stock Office colour scheme (dk1/lt1/dk2/lt2/
#      accent1-6/hlink/folHlink). Font scheme and format scheme (fills,
#      lines, effects) are identical between the two themes in this
#      file, so only the 12 colour-scheme slots actually change value.
# -----------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style on slide 16 (title + picture + one table -> table is
#    shape #3).
# ---------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$tableShape = $s16.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{C02D9BB8-2BE8-4244-AFEF-2932A9951E1F}")

# ---------------------------------------------------------------------
# 2) Theme colour scheme: Integral -> Office.
#    Colors(1..12) map to dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# ---------------------------------------------------------------------
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

$colorScheme.Colors(1).RGB  = RGB(0, 0, 0)          # dk1       000000
$colorScheme.Colors(2).RGB  = RGB(255, 255, 255)    # lt1       FFFFFF
$colorScheme.Colors(3).RGB  = RGB(68, 84, 106)      # dk2       44546A
$colorScheme.Colors(4).RGB  = RGB(231, 230, 230)    # lt2       E7E6E6
$colorScheme.Colors(5).RGB  = RGB(91, 155, 213)     # accent1   5B9BD5
$colorScheme.Colors(6).RGB  = RGB(237, 125, 49)     # accent2   ED7D31
$colorScheme.Colors(7).RGB  = RGB(165, 165, 165)    # accent3   A5A5A5
$colorScheme.Colors(8).RGB  = RGB(255, 192, 0)      # accent4   FFC000
$colorScheme.Colors(9).RGB  = RGB(68, 114, 196)     # accent5   4472C4
$colorScheme.Colors(10).RGB = RGB(112, 173, 71)     # accent6   70AD47
$colorScheme.Colors(11).RGB = RGB(5, 99, 193)       # hlink     0563C1
$colorScheme.Colors(12).RGB = RGB(149, 79, 114)     # folHlink  954F72
